$d = $word.ActiveDocument

$ldq = [char]0x201C
$rdq = [char]0x201D

# --- Edit 1: paragraph 2 ("Code is available in the src folder ... The screen") ---
# Collapses several runs (and the spell/gram proofErr markers around them) into a
# single run. Find/Replace across the run boundaries merges them and drops the
# now-interior proofErr markers.
$find1 = "Code is available in the src folder " + $ldq + "Searching.java" + $rdq + " . The screen"
$d.Content.Find.Execute($find1, $true, $false, $false, $false, $false, $true, 1, $false, $find1, 2) | Out-Null

# --- Edit 2: paragraph 7 ("The improvedLinearSearch () ... usual linear search.") ---
# Same run-merging/proofErr-removal, but also drops the trailing period and adds
# new trailing text (in two new runs) after it: ", when the array is uniformly
# distributed and sorted" and ".".
$find2 = "The improvedLinearSearch () takes in sorted array and compares the difference between key and starting and ending values of the array, after which it decides the end from which it will start iteration. Hence performs much better than usual linear search."
$replace2 = "The improvedLinearSearch () takes in sorted array and compares the difference between key and starting and ending values of the array, after which it decides the end from which it will start iteration. Hence performs much better than usual linear search, when the array is uniformly distributed and sorted."
$d.Content.Find.Execute($find2, $true, $false, $false, $false, $false, $true, 1, $false, $replace2, 2) | Out-Null

# Split the merged run so the new tail (", when the array is uniformly distributed
# and sorted" + ".") lives in its own two runs, matching the target markup.
$tailFind = ", when the array is uniformly distributed and sorted."
$hit = $d.Content.Find
$hit.Text = $tailFind
$hit.Forward = $true
$hit.Wrap = 1
$hit.Execute() | Out-Null
if ($hit.Found) {
    $tailStart = $d.Content.Start
}

$paraRange = $d.Paragraphs(7).Range
$full = $paraRange.Text
$markerOld = "usual linear search, when the array is uniformly distributed and sorted."
$idx = $full.IndexOf($markerOld)
if ($idx -ge 0) {
    $cut = $idx + ("usual linear search").Length
    $p1start = $paraRange.Start + $cut
    $p2start = $paraRange.Start + $full.IndexOf(", when the array is uniformly distributed and sorted.")
    $periodStart = $paraRange.Start + $full.LastIndexOf(".")

    $r2 = $d.Range($p2start, $periodStart)
    $r2.InsertParagraphAfter | Out-Null
}
